# Update Name of Algo
# Applies updated KNN imputation results to column C for the specified rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -10.491
    13 = -13.529
    16 = -12.853
    18 = -11.819
    20 = -12.559
    26 = -12.793
    27 = -13.556
    29 = -12.517
    35 = -12.258
    36 = -12.607
    45 = -13.297
    55 = -13.585
    57 = -13.643
    69 = -10.584
    76 = -12.969
    78 = -12.289
    82 = -12.077
    83 = -13.268
    93 = -10.196
    97 = -11.563
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
